$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (238) down into the
# new rows (239:244) so the new cells inherit the same styles
# (date format/border/alignment on column A, plain numbers on B:D).
$ws.Range("A238:D238").Copy() | Out-Null
$ws.Range("A239:D244").PasteSpecial(-4122) | Out-Null

# New data rows (dates 44313-44318 == 2021-04-27 .. 2021-05-02)
$data = @(
    @(44313, 4, 10, 117.827265229174),
    @(44314, 0, 10, 117.827265229174),
    @(44315, 1, 8, 94.26181218333922),
    @(44316, 1, 7, 82.47908566042182),
    @(44317, 2, 8, 94.26181218333922),
    @(44318, 4, 12, 141.3927182750088)
)

$r = 239
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value2 = $row[0]
    $ws.Cells.Item($r, 2).Value2 = $row[1]
    $ws.Cells.Item($r, 3).Value2 = $row[2]
    $ws.Cells.Item($r, 4).Value2 = $row[3]
    $r = $r + 1
}
